$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESC_Interface_All")

# Swap the "Value" column contents between row 8 (P2) and row 9 (P3)
$d8 = $ws.Range("D8").Value2
$d9 = $ws.Range("D9").Value2
$ws.Range("D8").Value2 = $d9
$ws.Range("D9").Value2 = $d8

# Update the selection to match the saved state (also clears the stale
# topLeftCell="A5" scroll position left over from the previous view)
$ws.Range("G16").Select()
